# "changed GUI to work with the new Main"
#
# The underlying ranking pipeline ("Main") was re-run, which re-derives the
# Twitter-influencer PageRank ranking. For this workbook that shows up as:
#   1) A handful of PageRank Score (column F) values for the top of the
#      table getting re-computed to a slightly different (but numerically
#      equivalent) floating point representation.
#   2) The bottom block of tied-score rows (rank 49 through rank 72, sheet
#      rows 50-54 and 57-73 - ranks 53/54 keep their original spot) being
#      re-ordered: same (User ID, Name, Username, Followers Count) tuples,
#      just shuffled into a different row order by the new script.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Score($row, $score) {
    $ws.Cells.Item($row, 6).Value = $score
}

function Set-Entry($row, $userId, $name, $username, $followers) {
    $ws.Cells.Item($row, 2).Value = $userId
    $ws.Cells.Item($row, 3).Value = $name
    $ws.Cells.Item($row, 4).Value = $username
    $ws.Cells.Item($row, 5).Value = $followers
}

# --- Re-derived PageRank Score values (column F), rows unaffected in data ---
Set-Score 3  0.006381806498112274
Set-Score 4  0.00332435198425819
Set-Score 5  0.0028982111442012984
Set-Score 6  0.0027244747627961285
Set-Score 7  0.0027071580345979025
Set-Score 15 0.0019842695136330505
Set-Score 17 0.0017786097712432381
Set-Score 20 0.0013437347789893891
Set-Score 21 0.0013085925479626491
Set-Score 22 0.0012412367405858636
Set-Score 25 0.0009736097081617733
Set-Score 26 0.0009020398400790572
Set-Score 27 0.0008953283935297851
Set-Score 29 0.0007696910001534317
Set-Score 33 0.0007124174743690127
Set-Score 38 0.0003570927401380605

# --- Re-ordered tail block (ranks 49-72), User ID / Name / Username / Followers ---
Set-Entry 50 940.0  "ChainGPT"            "@Chain_GPT"        1000000.0
Set-Entry 51 769.0  "Crypto Rover"        "@rovercrc"         980800.0
Set-Entry 52 978.0  "Everything Georgia"  "@GAFollowers"      1700000.0
Set-Entry 53 440.0  "POWR"                "@POWReSports"      702500.0
Set-Entry 54 853.0  "BSCDaily"            "@bsc_daily"        775100.0
Set-Entry 57 825.0  "Gordon"              "@AltcoinGordon"    539400.0
Set-Entry 58 268.0  "Revolving Games"     "@Revolving_Games"  592100.0
Set-Entry 59 331.0  "Crypto Rover"        "@rovercrc"         982300.0
Set-Entry 60 495.0  "Mnemonics_coin"      "@Mnemonics_coin"   568300.0
Set-Entry 61 212.0  "DuckChain"           "@Duck_Chain"       953900.0
Set-Entry 62 1022.0 "Chainlink"           "@chainlink"        1200000.0
Set-Entry 63 12.0   "Altcoin Daily"       "@AltcoinDailyio"   1700000.0
Set-Entry 64 845.0  "Somos Cosmos"        "@InformaCosmos"    1100000.0
Set-Entry 65 535.0  "Ice Open Network"    "@ice_blockchain"   2800000.0
Set-Entry 66 95.0   "Bitcoin.com News"    "@BTCTN"            3000000.0
Set-Entry 67 45.0   "Not Jerome Powell"   "@alifarhat79"      543000.0
Set-Entry 68 933.0  "Bitcoin Magazine"    "@BitcoinMagazine"  3400000.0
Set-Entry 69 772.0  "Gordon"              "@AltcoinGordon"    538500.0
Set-Entry 70 189.0  "Soompi"              "@soompi"           3600000.0
Set-Entry 71 1054.0 "Polyhedra"           "@PolyhedraZK"      965600.0
Set-Entry 72 340.0  "PepeMeme"            "@PepeMeme_"        536100.0
Set-Entry 73 852.0  "Ice Open Network"    "@ice_blockchain"   2800000.0
